$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: parse the order and store the per-item quantities in row 2 ---
$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = 4
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 3
$ws1.Range("E2").Value = 5
$ws1.Range("F2").Value = 2
$ws1.Range("G2").Value = 0
$ws1.Range("H2").Value = 0

# --- Check validity of the entered data: replace the old list-based
#     validations with a simple whole-number range check on the
#     "District Hospital" code cell ---
$ws1.Range("A2").Validation.Delete()
$ws1.Range("B2").Validation.Delete()
$ws1.Range("C2:H2").Validation.Delete()
$ws1.Range("A2").Validation.Add(1, 1, 1, 1, 2)

# --- Update the selected/active cells to reflect where data entry left off ---
$ws2.Range("B1:B49").Select()
$ws1.Activate()
$ws1.Range("H2").Select()
